# Update "想去人数" (interest count) figures in F column across sheets
# 展览 (Exhibitions), 演出 (Performances) and 全部类型 (All types) sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F5").Value = 969
$ws.Range("F6").Value = 363
$ws.Range("F8").Value = 556
$ws.Range("F9").Value = 1444
$ws.Range("F11").Value = 1335
$ws.Range("F12").Value = 2998
$ws.Range("F13").Value = 413
$ws.Range("F14").Value = 1609
$ws.Range("F16").Value = 791
$ws.Range("F17").Value = 236
$ws.Range("F18").Value = 1382
$ws.Range("F19").Value = 265
$ws.Range("F20").Value = 64
$ws.Range("F21").Value = 1120
$ws.Range("F22").Value = 397
$ws.Range("F23").Value = 3471
$ws.Range("F24").Value = 679
$ws.Range("F26").Value = 1533

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F7").Value = 50
$ws.Range("F13").Value = 15

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F11").Value = 50
$ws.Range("F15").Value = 969
$ws.Range("F16").Value = 363
$ws.Range("F18").Value = 556
$ws.Range("F19").Value = 1444
$ws.Range("F21").Value = 1335
$ws.Range("F22").Value = 2998
$ws.Range("F23").Value = 413
$ws.Range("F24").Value = 1609
$ws.Range("F26").Value = 791
$ws.Range("F27").Value = 236
$ws.Range("F28").Value = 1382
$ws.Range("F29").Value = 265
$ws.Range("F30").Value = 64
$ws.Range("F33").Value = 1120
$ws.Range("F34").Value = 397
$ws.Range("F35").Value = 3471
$ws.Range("F36").Value = 679
$ws.Range("F38").Value = 1533
$ws.Range("F40").Value = 15
